# Adapt the "asignacion_temas" form to add a new "Orden" (Order) column,
# and register it in the "referencia" help sheet.
$wb = $excel.ActiveWorkbook
$wsTemas = $wb.Worksheets.Item("temas")
$wsRef = $wb.Worksheets.Item("referencia")

# --- Sheet "referencia" ---
# Document the new "D" / "Orden" column with a new reference row.
$wsRef.Range("E10").Value = "D"
$wsRef.Range("F10").Value = "Orden"
$wsRef.Range("G10").Value = "Número que indica el orden en el que debe aparecer el tema dentro del programa"
$wsRef.Range("H10").Value = 1

# Extend the "description line" formula down into the new row (re-asserting
# I7:I10 with the same relative formula keeps I7:I9's existing shared group intact).
$wsRef.Range("I7:I10").FormulaR1C1 = '=RC[-4] & ", " & RC[-3] & "," & RC[-2] & "," & RC[-1]'

$wsRef.Range("E10:I10").Select()

# --- Sheet "temas" ---
# Header D1: clone the header formatting (bold / filled) from C1, then set its text.
$wsTemas.Range("C1").Copy()
$wsTemas.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
$wsTemas.Range("D1").Value = "Orden"

# Data D2:D4: clone the plain data-cell formatting from C2, then fill in the order numbers.
$wsTemas.Range("C2").Copy()
$wsTemas.Range("D2").PasteSpecial(-4104)  # xlPasteAll
$wsTemas.Range("D2").Value = 1

$wsTemas.Range("C2").Copy()
$wsTemas.Range("D3").PasteSpecial(-4104)  # xlPasteAll
$wsTemas.Range("D3").Value = 2

$wsTemas.Range("C2").Copy()
$wsTemas.Range("D4").PasteSpecial(-4104)  # xlPasteAll
$wsTemas.Range("D4").Value = 3

$excel.CutCopyMode = 0

# Move the selection to the first unfrozen row below the data, matching the
# author's final cursor position, and leave "temas" as the active sheet/tab.
$wsTemas.Activate()
$wsTemas.Range("A5:XFD5").Select()

$wb.Save()
